$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pdfgen")

# Insert a new row above the old row 5 ("skeleton /fonts/openfonts"), pushing
# it (and the rows below) down by one.
$ws.Rows.Item(5).Insert()

# Copy the formatting of the row above (row 4) onto the newly inserted blank
# row so it matches the rest of the table.
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)

# Row 4: the /bin/fop skeleton copy becomes a "private" entry (hidden from UI
# with the leading '#') since it is now only used as the source for compare.
$ws.Range("A4").Value = "#skeleton"

# New row 5: compare the freshly-deployed /bin/fop against the skeleton copy
# and report any differences, so pdfgen's fop install stays in sync.
$ws.Range("A5").Value = "compare"
$ws.Range("B5").Value = "/bin/fop"
$ws.Range("C5").Value = "skeleton"
$ws.Range("D5").Value = "report"

$ws.Range("D16").Select()
